$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row (TCID 48 - Leave Application Over-Utilization, Hourly)
$ws.Range("A43").Value = "48"
$ws.Range("B43").Value = "LeaveApplicaton"
$ws.Range("C43").Value = "OverUtilization_Hourly"
$ws.Range("D43").Value = "com.darwinbox.leaves.Application.OverUtilization_Hourly"
$ws.Range("E43").Value = "Application//LeaveApplication.xlsx"
$ws.Range("F43").Value = "OverUtilizationScenarioesHourly"
$ws.Range("G43").Value = "All"

# Widen column F so the new SheetName text fits (best-fit sizing)
$ws.Columns("F").ColumnWidth = 30.109375

# Reflect the cursor/selection ending on the newly added row
$ws.Range("F43").Select() | Out-Null
